# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Reverses the order of the "Periodo Mora" (period) values shown in rows
# 16..32 (newest period first), and refreshes the "Valor Mora" (G column)
# amounts.  The "Salario Basico" (F column) value of 25740 stays attached
# to period 1903 (now on the first data row) while every other period
# keeps 27578; all Valor Mora values become 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1903","1902","1901","1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801","1709","1708")

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value = $period   # column E - Periodo Mora

    if ($period -eq "1903") {
        $ws.Cells.Item($row, 6).Value = 25740 # column F - Salario Basico
    } else {
        $ws.Cells.Item($row, 6).Value = 27578
    }

    $ws.Cells.Item($row, 7).Value = 781242    # column G - Valor Mora
}
